# Update "Last Updated" timestamp on the Metadata sheet.
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 12:38 PM"

# The "Stock List" sheet refreshed: every row's data (Stock/Stock Name/
# Price/% Change/Market Cap) shifted up by one - row N now holds what
# row N+1 held before - and a brand new row was appended at the bottom
# (row 76) with fresh data (TRAVELFOOD).
$ws = $wb.Worksheets.Item("Stock List")

for ($r = 2; $r -le 75; $r++) {
    $nr = $r + 1
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($nr, 2).Value2
    $ws.Cells.Item($r, 3).Value = $ws.Cells.Item($nr, 3).Value2
    $ws.Cells.Item($r, 4).Value = $ws.Cells.Item($nr, 4).Value2
    $ws.Cells.Item($r, 5).Value = $ws.Cells.Item($nr, 5).Value2
    $ws.Cells.Item($r, 8).Value = $ws.Cells.Item($nr, 8).Value2
}

$ws.Cells.Item(76, 2).Value = "TRAVELFOOD"
$ws.Cells.Item(76, 3).Value = "TRAVELFOOD"
$ws.Cells.Item(76, 4).Value = 1316.3
$ws.Cells.Item(76, 5).Value = 0.1141
$ws.Cells.Item(76, 8).Value = 17332.9705
